$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.165747
$ws.Range("H2").Value = 0.497241
$ws.Range("I2").Value = 0.008095785894995438
$ws.Range("J2").Value = 0.00809578589499544
$ws.Range("M2").Value = 2.759544333333333
$ws.Range("N2").Value = 8.278632999999999
$ws.Range("O2").Value = 0.2574067337278401
$ws.Range("P2").Value = 0.2574067337278401
$ws.Range("Q2").Value = 0.457386194617
$ws.Range("R2").Value = 4.116475751553
$ws.Range("S2").Value = 0.002083909804190694
$ws.Range("T2").Value = 0.002083909804190695
$ws.Range("G3").Value = 0.165747
$ws.Range("H3").Value = 0.497241
$ws.Range("I3").Value = 0.008095785894995438
$ws.Range("J3").Value = 0.00809578589499544
$ws.Range("O3").Value = 0.6758254232987829
$ws.Range("P3").Value = 0.6758254232987829
$ws.Range("Q3").Value = 1.200874639569
$ws.Range("R3").Value = 10.807871756121
$ws.Range("S3").Value = 0.005471337929421608
$ws.Range("T3").Value = 0.00547133792942161
$ws.Range("G4").Value = 0.165747
$ws.Range("H4").Value = 0.497241
$ws.Range("I4").Value = 0.008095785894995438
$ws.Range("J4").Value = 0.00809578589499544
$ws.Range("M4").Value = 0.5200313333333334
$ws.Range("N4").Value = 1.560094
$ws.Range("O4").Value = 0.0485078515798926
$ws.Range("P4").Value = 0.0485078515798926
$ws.Range("Q4").Value = 0.08619363340600002
$ws.Range("R4").Value = 0.7757427006540001
$ws.Range("S4").Value = 0.0003927091806170267
$ws.Range("T4").Value = 0.0003927091806170267
$ws.Range("G5").Value = 0.165747
$ws.Range("H5").Value = 0.497241
$ws.Range("I5").Value = 0.008095785894995438
$ws.Range("J5").Value = 0.00809578589499544
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1957573333333333
$ws.Range("N5").Value = 0.587272
$ws.Range("O5").Value = 0.01825999139348442
$ws.Range("P5").Value = 0.01825999139348442
$ws.Range("Q5").Value = 0.032446190728
$ws.Range("R5").Value = 0.292015716552
$ws.Range("S5").Value = 0.0001478289807661093
$ws.Range("T5").Value = 0.0001478289807661093
$ws.Range("I6").Value = 0.7079722685862583
$ws.Range("J6").Value = 0.7079722685862583
$ws.Range("M6").Value = 2.759544333333333
$ws.Range("N6").Value = 8.278632999999999
$ws.Range("O6").Value = 0.2574067337278401
$ws.Range("P6").Value = 0.2574067337278401
$ws.Range("Q6").Value = 39.99818498451233
$ws.Range("R6").Value = 359.983664860611
$ws.Range("S6").Value = 0.1822368292266779
$ws.Range("T6").Value = 0.1822368292266779
$ws.Range("I7").Value = 0.7079722685862583
$ws.Range("J7").Value = 0.7079722685862583
$ws.Range("O7").Value = 0.6758254232987829
$ws.Range("P7").Value = 0.6758254232987829
$ws.Range("S7").Value = 0.4784656581011076
$ws.Range("T7").Value = 0.4784656581011076
$ws.Range("I8").Value = 0.7079722685862583
$ws.Range("J8").Value = 0.7079722685862583
$ws.Range("M8").Value = 0.5200313333333334
$ws.Range("N8").Value = 1.560094
$ws.Range("O8").Value = 0.0485078515798926
$ws.Range("P8").Value = 0.0485078515798926
$ws.Range("Q8").Value = 7.537588440655335
$ws.Range("R8").Value = 67.83829596589801
$ws.Range("S8").Value = 0.03434221372726208
$ws.Range("T8").Value = 0.03434221372726207
$ws.Range("I9").Value = 0.7079722685862583
$ws.Range("J9").Value = 0.7079722685862583
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1957573333333333
$ws.Range("N9").Value = 0.587272
$ws.Range("O9").Value = 0.01825999139348442
$ws.Range("P9").Value = 0.01825999139348442
$ws.Range("Q9").Value = 2.837402514669334
$ws.Range("R9").Value = 25.536622632024
$ws.Range("S9").Value = 0.01292756753121072
$ws.Range("T9").Value = 0.01292756753121072
$ws.Range("G10").Value = 5.642879333333333
$ws.Range("H10").Value = 16.928638
$ws.Range("I10").Value = 0.2756221404547972
$ws.Range("J10").Value = 0.2756221404547972
$ws.Range("M10").Value = 2.759544333333333
$ws.Range("N10").Value = 8.278632999999999
$ws.Range("O10").Value = 0.2574067337278401
$ws.Range("P10").Value = 0.2574067337278401
$ws.Range("Q10").Value = 15.57177568798378
$ws.Range("R10").Value = 140.145981191854
$ws.Range("S10").Value = 0.0709469949175453
$ws.Range("T10").Value = 0.0709469949175453
$ws.Range("G11").Value = 5.642879333333333
$ws.Range("H11").Value = 16.928638
$ws.Range("I11").Value = 0.2756221404547972
$ws.Range("J11").Value = 0.2756221404547972
$ws.Range("O11").Value = 0.6758254232987829
$ws.Range("P11").Value = 0.6758254232987829
$ws.Range("Q11").Value = 40.88394170360866
$ws.Range("R11").Value = 367.955475332478
$ws.Range("S11").Value = 0.1862724497433799
$ws.Range("T11").Value = 0.1862724497433799
$ws.Range("G12").Value = 5.642879333333333
$ws.Range("H12").Value = 16.928638
$ws.Range("I12").Value = 0.2756221404547972
$ws.Range("J12").Value = 0.2756221404547972
$ws.Range("M12").Value = 0.5200313333333334
$ws.Range("N12").Value = 1.560094
$ws.Range("O12").Value = 0.0485078515798926
$ws.Range("P12").Value = 0.0485078515798926
$ws.Range("Q12").Value = 2.934474063552445
$ws.Range("R12").Value = 26.410266571972
$ws.Range("S12").Value = 0.01336983788131361
$ws.Range("T12").Value = 0.01336983788131361
$ws.Range("G13").Value = 5.642879333333333
$ws.Range("H13").Value = 16.928638
$ws.Range("I13").Value = 0.2756221404547972
$ws.Range("J13").Value = 0.2756221404547972
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1957573333333333
$ws.Range("N13").Value = 0.587272
$ws.Range("O13").Value = 0.01825999139348442
$ws.Range("P13").Value = 0.01825999139348442
$ws.Range("Q13").Value = 1.104635010615111
$ws.Range("R13").Value = 9.941715095536001
$ws.Range("S13").Value = 0.005032857912558351
$ws.Range("T13").Value = 0.005032857912558351
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.1701286666666667
$ws.Range("H14").Value = 0.510386
$ws.Range("I14").Value = 0.008309805063949155
$ws.Range("J14").Value = 0.008309805063949155
$ws.Range("M14").Value = 2.759544333333333
$ws.Range("N14").Value = 8.278632999999999
$ws.Range("O14").Value = 0.2574067337278401
$ws.Range("P14").Value = 0.2574067337278401
$ws.Range("Q14").Value = 0.4694775980375555
$ws.Range("R14").Value = 4.225298382338
$ws.Range("S14").Value = 0.002138999779426217
$ws.Range("T14").Value = 0.002138999779426217
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.1701286666666667
$ws.Range("H15").Value = 0.510386
$ws.Range("I15").Value = 0.008309805063949155
$ws.Range("J15").Value = 0.008309805063949155
$ws.Range("O15").Value = 0.6758254232987829
$ws.Range("P15").Value = 0.6758254232987829
$ws.Range("Q15").Value = 1.232620809207333
$ws.Range("R15").Value = 11.093587282866
$ws.Range("S15").Value = 0.005615977524873808
$ws.Range("T15").Value = 0.005615977524873808
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.1701286666666667
$ws.Range("H16").Value = 0.510386
$ws.Range("I16").Value = 0.008309805063949155
$ws.Range("J16").Value = 0.008309805063949155
$ws.Range("M16").Value = 0.5200313333333334
$ws.Range("N16").Value = 1.560094
$ws.Range("O16").Value = 0.0485078515798926
$ws.Range("P16").Value = 0.0485078515798926
$ws.Range("Q16").Value = 0.08847223736488891
$ws.Range("R16").Value = 0.7962501362840001
$ws.Range("S16").Value = 0.0004030907906998856
$ws.Range("T16").Value = 0.0004030907906998855
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.1701286666666667
$ws.Range("H17").Value = 0.510386
$ws.Range("I17").Value = 0.008309805063949155
$ws.Range("J17").Value = 0.008309805063949155
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1957573333333333
$ws.Range("N17").Value = 0.587272
$ws.Range("O17").Value = 0.01825999139348442
$ws.Range("P17").Value = 0.01825999139348442
$ws.Range("Q17").Value = 0.03330393411022223
$ws.Range("R17").Value = 0.299735406992
$ws.Range("S17").Value = 0.0001517369689492448
$ws.Range("T17").Value = 0.0001517369689492448
